# Update cryptocurrency price/volume figures (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "29.304.55"
$ws.Range("D3").Value = "1.876.45"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7118"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3103"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07768"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.21%  "
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08510"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.09%  "
$ws.Range("D12").Value = "1.879.63"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7098"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").Value = "29.306.58"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008247"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.004"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "2.133.52"
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.817"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1622"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.021"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E28").Value = "  +1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.511"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.401"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.320"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.278"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05236"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.930"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.177"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7397"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.687"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01864"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.36%  "
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("D40").Value = "1.174.98"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("E41").Value = "  +3.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8890"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "2.029.65"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("E47").Value = "  +2.49%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000122"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.383"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4306"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.95%  "
